$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing date value in D2 (Mon Jun 14 2021 -> Fri Jun 18 2021)
$ws.Range("D2").Value = "Fri Jun 18 2021"

# Append a new data row (row 3) mirroring the structure of row 2
$ws.Range("A3").Value = "bar 2"
$ws.Range("B3").Value = "moo 2"
$ws.Range("C3").Value = 1234
$ws.Range("D3").Value = "Wed May 19 2021"
